$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows just above the old row 27 ("27-73" / A0 / A1 / A2 / A3 /
# "58 (A4)" / "59 (A5)" block), pushing that block (old rows 27-33) down to 29-35.
# xlShiftDown = -4121
$null = $ws.Range("A27:B28").Insert(-4121)

# New DIO 12 / 13 are now wired to the Servo.
$ws.Range("B12").Value = "Servo"
$ws.Range("B13").Value = "Servo"

# New DIO 52 / 53 rows for the wait button (order of writes controls the
# shared-string table order: HIGH ends up before "Wait button").
$ws.Range("B28").Value = "HIGH"
$ws.Range("A27").Value = 52
$ws.Range("B27").Value = "Wait button"
$ws.Range("A28").Value = 53

# Restore the sheet view: scrolled so row 4 is at the top, with B14 selected.
$null = $ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$null = $ws.Range("B14").Select()
